# "fix for experiment 1"
# The K=2 block of trials (rows 19-30) recorded the wrong winner ("X" vs "V")
# on four of the twelve trials. Correct the outcome cells; the COUNTIF /
# win-chance formulas in rows 31-33 (and the rolled-up summary in rows 53-55)
# recalculate automatically from these corrected raw results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "V"
$ws.Range("E24").Value = "X"

$ws.Range("B28").Value = "V"
$ws.Range("E28").Value = "X"

$ws.Range("C29").Value = "X"
$ws.Range("E29").Value = "V"

$ws.Range("D30").Value = "V"
$ws.Range("E30").Value = "X"

# Reflect the scrolled/zoomed view the author was looking at while making the fix.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F21").Select()
